$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertAfter("Test sentence 3.")

$lastPara2 = $d.Paragraphs.Last
$r2 = $lastPara2.Range
$r2.InsertParagraphAfter()
$r2.Collapse(0)
$r2.InsertAfter("Test sentence 4.")
